$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new rows 12-17 in columns A:D, mirroring the abbreviated
# "NA" sample-group rows already present in columns K:N (rows 4-9),
# extending the full Sample group / Replicate / Experimental Group / CustomID
# table that begins at row 3.

$data = @(
    @(12, "NA", 1, "Treatment", "NA_T_1"),
    @(13, "NA", 2, "Treatment", "NA_T_2"),
    @(14, "NA", 3, "Treatment", "NA_T_3"),
    @(15, "NA", 1, "Control",   "NA_1"),
    @(16, "NA", 2, "Control",   "NA_2"),
    @(17, "NA", 3, "Control",   "NA_3")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("A$r").Value = $row[1]
    $ws.Range("B$r").Value = $row[2]
    $ws.Range("C$r").Value = $row[3]
    $ws.Range("D$r").Value = $row[4]
}

# Update the view: selection and scroll position
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("K4:N15").Select() | Out-Null
